$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at position 1256; this shifts the existing rows
# 1256..1358 down to 1257..1359, preserving their original values/format.
$ws.Rows("1256:1256").Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A1256").Value = 6
$ws.Range("B1256").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1256").Value = "Metropolitana"
$ws.Range("D1256").Value = 45166
$ws.Range("E1256").Value = 13
$ws.Range("F1256").Value = 100112031
$ws.Range("G1256").Value = "Poroto verde"
$ws.Range("H1256").Value = "Magnum"
$ws.Range("I1256").Value = "Primera"
$ws.Range("J1256").Value = 250
$ws.Range("K1256").Value = 20000
$ws.Range("L1256").Value = 25000
$ws.Range("M1256").Value = 22000
$ws.Range("N1256").Value = '$/malla 25 kilos'
$ws.Range("O1256").Value = "Perú"
$ws.Range("P1256").Value = 880
$ws.Range("Q1256").Value = 25
$ws.Range("R1256").Value = "Hortaliza"
